$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "urban"/"rural" (шаар/айыл, город/село) rows with the new
# "City"/"Village" (Шаар жерлери/Айыл аймагы, Городские поселения/Сельская
# местность) rows. This also causes the now-unused shared strings
# ("город", "село", "urban", "rural") to be dropped from the workbook.
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Update the Kyrgyz sub-title text (done last so it lands at the end of the
# shared-string table, matching the target ordering).
$ws.Range("A2").Value = "(жалпы калктын санына карата пайыз менен)"

# Update the saved selection/scroll position of the sheet view.
$ws.Range("A8").Select()
